# Atualizei dados da bibi e add
# Insert two new daily-revenue rows for August (dias 12 e 13) right after the
# existing August block (rows 2-12, dias 1-11). Inserting physical rows
# pushes every subsequent row (July/June/May data) down by two, which is
# exactly the shift seen across the rest of the sheet in the diff - no other
# cell needs to be touched by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 13:14 - everything from row 13 down shifts by +2.
$ws.Rows("13:14").Insert()

# Fill the newly inserted rows with the new daily totals for August/2025.
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 28398.08
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 2025
$ws.Range("E13").Value = "08/2025"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 25151.7
$ws.Range("C14").Value = 8
$ws.Range("D14").Value = 2025
$ws.Range("E14").Value = "08/2025"
